$d = $word.ActiveDocument

$lsq = [char]8216   # U+2018 LEFT SINGLE QUOTATION MARK
$rsq = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK

$findText = "called " + $lsq + "Mapping" + $rsq
$replaceText = "called " + $lsq + "mapping" + $rsq

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
